$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / label / URL updates (Coin, Link, Volume columns)
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E18').Value = '17OneONEWorstin24h'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E21').Value = '20BitKanKAN'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('B23').Value = 'UpBots'
$ws.Range('C23').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('E23').Value = '22UpBotsUBXT'
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('E43').Value = '42CEJICEJI'
$ws.Range('E47').Value = '46ACDXExchangeACXT'
$ws.Range('E49').Value = '48BOLOBOLOBestin24h'

# Numeric-looking price values stored as literal text in the source data.
# Preserve their Text cell type (avoid Excel auto-converting to a Number)
# by temporarily formatting as Text, then restoring the original style.
$cell = $ws.Range('D2')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '244.36'
$cell.Style = $origStyle

$cell = $ws.Range('D3')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '21.49'
$cell.Style = $origStyle

$cell = $ws.Range('D4')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.222'
$cell.Style = $origStyle

$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.05605'
$cell.Style = $origStyle

$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.364'
$cell.Style = $origStyle

$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.373'
$cell.Style = $origStyle

$cell = $ws.Range('D8')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.8050'
$cell.Style = $origStyle

$cell = $ws.Range('D9')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.9440'
$cell.Style = $origStyle

$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.07304'
$cell.Style = $origStyle

$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.03148'
$cell.Style = $origStyle

$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.03055'
$cell.Style = $origStyle

$cell = $ws.Range('D14')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.09282'
$cell.Style = $origStyle

$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.560'
$cell.Style = $origStyle

$cell = $ws.Range('D16')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.001664'
$cell.Style = $origStyle

$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.04705'
$cell.Style = $origStyle

$cell = $ws.Range('D18')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0005746'
$cell.Style = $origStyle

$cell = $ws.Range('D19')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.006427'
$cell.Style = $origStyle

$cell = $ws.Range('D20')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.004991'
$cell.Style = $origStyle

$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.001042'
$cell.Style = $origStyle

$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0001498'
$cell.Style = $origStyle

$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0003098'
$cell.Style = $origStyle

$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.758'
$cell.Style = $origStyle

$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.094'
$cell.Style = $origStyle

$cell = $ws.Range('D26')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.3261'
$cell.Style = $origStyle

$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.03909'
$cell.Style = $origStyle

$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.006891'
$cell.Style = $origStyle

$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.1037'
$cell.Style = $origStyle

$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.002896'
$cell.Style = $origStyle

$cell = $ws.Range('D44')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.008501'
$cell.Style = $origStyle

$cell = $ws.Range('D45')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.00005926'
$cell.Style = $origStyle

$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.00000000749'
$cell.Style = $origStyle

$cell = $ws.Range('D47')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0005496'
$cell.Style = $origStyle

$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.6820'
$cell.Style = $origStyle

$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.08442'
$cell.Style = $origStyle

$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.01009'
$cell.Style = $origStyle
